$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header (row 2), pushing every existing
# member row down by one. Every pre-existing member keeps exactly the
# data it had before - only its row number shifts by one.
$ws.Rows.Item(2).Insert()

# New member row: "teste update dois" with "-" placeholders in every
# category column (same convention already used for every other empty
# cell in the sheet).
$ws.Range("A2").Value = "teste update dois"
$ws.Range("B2:E2").Value = "-"

# Row 3 (the old row 2, just shifted down) still carries the correct
# plain body font plus the "blue" zebra-stripe fill that the new row
# should use too. Insert() left row 2 with the header's bold/white font,
# so clone row 3's whole format onto row 2 to turn it into a normal body
# row before fixing up the banding below.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Re-apply the report's zebra-stripe row banding (even data row -> the
# "blue" fill, odd data row -> the "white" fill) across every data row so
# the stripes line up with the new row positions.
for ($r = 2; $r -le 9; $r++) {
    $row = $ws.Range("A" + $r + ":E" + $r)
    if (($r % 2) -eq 0) {
        $row.Interior.ColorIndex = 35
    } else {
        $row.Interior.ColorIndex = 2
    }
}
